$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new columns B (First release) and C (Second release) ---
$ws.Range("B1").Value = "First release"
$ws.Range("C1").Value = "Second release"

# Copy the header formatting (bold maroon on grey fill) from A1 onto the
# new header cells so they reuse the existing "header" style instead of
# Excel fabricating a brand new style entry.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1:C1").PasteSpecial(-4122) | Out-Null

# --- Mark protocols published in the first release (column B) ---
$ws.Range("B2").Value = "p"
$ws.Range("B3").Value = "p"
$ws.Range("B4").Value = "p"
$ws.Range("B7").Value = "p"
$ws.Range("B8").Value = "p"
$ws.Range("B9").Value = "p"

# --- Mark protocols published in the second release (column C) ---
$ws.Range("C17").Value = "p"
$ws.Range("C18").Value = "p"
$ws.Range("C21").Value = "p"
$ws.Range("C22").Value = "p"
$ws.Range("C23").Value = "p"
$ws.Range("C27").Value = "p"

# --- Column widths for the two new columns ---
$ws.Columns.Item(2).ColumnWidth = 19.25
$ws.Columns.Item(3).ColumnWidth = 13.42

# --- Selection left where the author's last edit landed ---
$ws.Range("B51").Select() | Out-Null

Write-Output "done"
